# "Error Calculations and Plots" - refresh the imputed/missing-data sample:
#  - two ID rows (RM 232, SC 92) are dropped from the table entirely, which
#    shrinks the used range from A1:F35 down to A1:F33 and shifts every row
#    below them up
#  - a handful of individual cells in columns C/E/F are swapped between
#    "missing" (blank) and a concrete imputed number

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 (ID "RM 232") is removed outright.
$ws.Rows.Item(26).Delete()

# After that shift, the old row 28 (ID "SC 92") is now row 27 - remove it too.
$ws.Rows.Item(27).Delete()

# Remaining per-cell value changes, addressed by their FINAL (post-delete) row
# numbers now that the sheet only spans A1:F33.
$ws.Range("E2").Value = -7.2
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("F8").Value = 17.05
$ws.Range("F10").Value = 16.43
$ws.Range("E11").Value = -7.9
$ws.Range("F12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("F15").Value = 16.2
$ws.Range("F18").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E25").ClearContents()
$ws.Range("F25").Value = 16.6
$ws.Range("F27").Value = 17
$ws.Range("C29").ClearContents()
$ws.Range("F29").ClearContents()
$ws.Range("C33").Value = 10.4
$ws.Range("E33").Value = -10.7
$ws.Range("F33").ClearContents()
